$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.554972
$ws.Range("H2").Value = 10.664916
$ws.Range("I2").Value = 0.01373869946676525
$ws.Range("J2").Value = 0.01373869946676525
$ws.Range("M2").Value = 1.438916666666667
$ws.Range("N2").Value = 4.31675
$ws.Range("O2").Value = 0.07110575088422412
$ws.Range("P2").Value = 0.07110575088422412
$ws.Range("Q2").Value = 5.115308460333333
$ws.Range("R2").Value = 46.03777614299999
$ws.Range("S2").Value = 0.0009769005417570324
$ws.Range("T2").Value = 0.0009769005417570326
$ws.Range("G3").Value = 3.554972
$ws.Range("H3").Value = 10.664916
$ws.Range("I3").Value = 0.01373869946676525
$ws.Range("J3").Value = 0.01373869946676525
$ws.Range("O3").Value = 0.07524532698869241
$ws.Range("P3").Value = 0.07524532698869241
$ws.Range("Q3").Value = 5.413107279782666
$ws.Range("R3").Value = 48.717965518044
$ws.Range("S3").Value = 0.001033772933776125
$ws.Range("T3").Value = 0.001033772933776125
$ws.Range("G4").Value = 3.554972
$ws.Range("H4").Value = 10.664916
$ws.Range("I4").Value = 0.01373869946676525
$ws.Range("J4").Value = 0.01373869946676525
$ws.Range("N4").Value = 51.824064
$ws.Range("O4").Value = 0.8536489221270834
$ws.Range("P4").Value = 0.8536489221270835
$ws.Range("Q4").Value = 61.411032148736
$ws.Range("R4").Value = 552.699289338624
$ws.Range("S4").Value = 0.01172802599123209
$ws.Range("T4").Value = 0.01172802599123209
$ws.Range("I5").Value = 0.9249383070209767
$ws.Range("J5").Value = 0.9249383070209768
$ws.Range("M5").Value = 1.438916666666667
$ws.Range("N5").Value = 4.31675
$ws.Range("O5").Value = 0.07110575088422412
$ws.Range("P5").Value = 0.07110575088422412
$ws.Range("Q5").Value = 344.3808315798889
$ws.Range("R5").Value = 3099.427484219
$ws.Range("S5").Value = 0.06576843284230957
$ws.Range("T5").Value = 0.06576843284230958
$ws.Range("I6").Value = 0.9249383070209767
$ws.Range("J6").Value = 0.9249383070209768
$ws.Range("O6").Value = 0.07524532698869241
$ws.Range("P6").Value = 0.07524532698869241
$ws.Range("S6").Value = 0.06959728535616096
$ws.Range("T6").Value = 0.06959728535616097
$ws.Range("I7").Value = 0.9249383070209767
$ws.Range("J7").Value = 0.9249383070209768
$ws.Range("N7").Value = 51.824064
$ws.Range("O7").Value = 0.8536489221270834
$ws.Range("P7").Value = 0.8536489221270835
$ws.Range("S7").Value = 0.7895725888225061
$ws.Range("T7").Value = 0.7895725888225064
$ws.Range("I8").Value = 0.06132299351225805
$ws.Range("J8").Value = 0.06132299351225806
$ws.Range("M8").Value = 1.438916666666667
$ws.Range("N8").Value = 4.31675
$ws.Range("O8").Value = 0.07110575088422412
$ws.Range("P8").Value = 0.07110575088422412
$ws.Range("Q8").Value = 22.83229415455556
$ws.Range("R8").Value = 205.490647391
$ws.Range("S8").Value = 0.004360417500157513
$ws.Range("T8").Value = 0.004360417500157514
$ws.Range("I9").Value = 0.06132299351225805
$ws.Range("J9").Value = 0.06132299351225806
$ws.Range("O9").Value = 0.07524532698869241
$ws.Range("P9").Value = 0.07524532698869241
$ws.Range("S9").Value = 0.00461426869875532
$ws.Range("T9").Value = 0.004614268698755321
$ws.Range("I10").Value = 0.06132299351225805
$ws.Range("J10").Value = 0.06132299351225806
$ws.Range("N10").Value = 51.824064
$ws.Range("O10").Value = 0.8536489221270834
$ws.Range("P10").Value = 0.8536489221270835
$ws.Range("S10").Value = 0.05234830731334522
$ws.Range("T10").Value = 0.05234830731334523
